# Fix acceltic format skill
# -------------------------------------------------------------
# Rebuilds the "Duplicates" sheet (Table13) so it matches the refreshed
# Skill/Parent ontology: 15 new rows of skill aliases are inserted, the whole
# A1:B66 block becomes A1:B81, and the table/autofilter/selection are resized
# to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Duplicates"
$ws.Activate()

# --- 1. Make room: insert 15 blank rows right before the existing block ---
$ws.Rows.Item(31).Resize(15).Insert()

# --- 2. Write the refreshed Skill (A) / Parent (B) pairs for rows 31-81 ---
$ws.Range("A31").Value = "AI"
$ws.Range("B31").Value = "Artificial Intelligence (AI)"
$ws.Range("A32").Value = "Athena"
$ws.Range("B32").Value = "Athena Cloud"
$ws.Range("A33").Value = "Auditing"
$ws.Range("B33").Value = "Audit"
$ws.Range("A34").Value = "Maya"
$ws.Range("B34").Value = "Autodesk Maya"
$ws.Range("A35").Value = "Blog"
$ws.Range("B35").Value = "Blogging"
$ws.Range("A36").Value = "Brainstorm"
$ws.Range("B36").Value = "Brainstorming"
$ws.Range("A37").Value = "Budgeting"
$ws.Range("B37").Value = "Budget"
$ws.Range("A38").Value = "C Sharp"
$ws.Range("B38").Value = "C#"
$ws.Range("A39").Value = "Cloud Technologies"
$ws.Range("B39").Value = "Cloud Technology"
$ws.Range("A40").Value = "Google Vision"
$ws.Range("B40").Value = "Cloud Vision API"
$ws.Range("A41").Value = "Google Cloud's Vision API"
$ws.Range("B41").Value = "Cloud Vision API"
$ws.Range("A42").Value = "Google Cloud Vision"
$ws.Range("B42").Value = "Cloud Vision API"
$ws.Range("A43").Value = "Hue"
$ws.Range("B43").Value = "Cloudera Hue"
$ws.Range("A44").Value = "D3"
$ws.Range("B44").Value = "D3.js"
$ws.Range("A45").Value = "Dashboards"
$ws.Range("B45").Value = "Dashboard"
$ws.Range("A46").Value = "Data Cleansing/Preparation"
$ws.Range("B46").Value = "Data Cleansing"
$ws.Range("A47").Value = "Data Storytelling & Presentation Skill"
$ws.Range("B47").Value = "Data Storytelling"
$ws.Range("A48").Value = "DBMS"
$ws.Range("B48").Value = "Database Management System"
$ws.Range("A49").Value = "ETL"
$ws.Range("B49").Value = "Extract Transform Load (ETL)"
$ws.Range("A50").Value = "Financial Statements"
$ws.Range("B50").Value = "Financial Statement"
$ws.Range("A51").Value = "Go"
$ws.Range("B51").Value = "Golang"
$ws.Range("A52").Value = "Google Adwords"
$ws.Range("B52").Value = "Google Ads"
$ws.Range("A53").Value = "Adwords"
$ws.Range("B53").Value = "Google Ads"
$ws.Range("A54").Value = "Google Cloud"
$ws.Range("B54").Value = "Google Cloud Platform (GCP)"
$ws.Range("A55").Value = "G Suite"
$ws.Range("B55").Value = "Google Workspace"
$ws.Range("A56").Value = "IntelliJ"
$ws.Range("B56").Value = "IntelliJ IDEA"
$ws.Range("A57").Value = "k-Nearest Neighbour (k-NN)"
$ws.Range("B57").Value = "k-Nearest Neighbour (kNN)"
$ws.Range("A58").Value = "Azure"
$ws.Range("B58").Value = "Microsoft Azure"
$ws.Range("A59").Value = "Excel"
$ws.Range("B59").Value = "Microsoft Excel"
$ws.Range("A60").Value = "MSExcel"
$ws.Range("B60").Value = "Microsoft Excel"
$ws.Range("A61").Value = "OneNote"
$ws.Range("B61").Value = "Microsoft OneNote"
$ws.Range("A62").Value = "Outlook"
$ws.Range("B62").Value = "Microsoft Outlook"
$ws.Range("A63").Value = "Power BI"
$ws.Range("B63").Value = "Microsoft Power BI"
$ws.Range("A64").Value = "PowerBI"
$ws.Range("B64").Value = "Microsoft Power BI"
$ws.Range("A65").Value = "PowerPoint"
$ws.Range("B65").Value = "Microsoft PowerPoint"
$ws.Range("A66").Value = "MSWord"
$ws.Range("B66").Value = "Microsoft Word"
$ws.Range("A67").Value = "Mongo DB"
$ws.Range("B67").Value = "MongoDB"
$ws.Range("A68").Value = "NLP"
$ws.Range("B68").Value = "Natural Language Processing (NLP)"
$ws.Range("A69").Value = "Node"
$ws.Range("B69").Value = "Node.js"
$ws.Range("A70").Value = "Python Programming"
$ws.Range("B70").Value = "Python"
$ws.Range("A71").Value = "R Programming"
$ws.Range("B71").Value = "R"
$ws.Range("A72").Value = "React.js"
$ws.Range("B72").Value = "React"
$ws.Range("A73").Value = "bazingaJS"
$ws.Range("B73").Value = "React"
$ws.Range("A74").Value = "Regressions"
$ws.Range("B74").Value = "Regression"
$ws.Range("A75").Value = "Reinforcements"
$ws.Range("B75").Value = "Reinforcement Learning"
$ws.Range("A76").Value = "Rails"
$ws.Range("B76").Value = "Ruby on Rails"
$ws.Range("A77").Value = "BusinessObjects"
$ws.Range("B77").Value = "SAP BusinessObjects"
$ws.Range("A78").Value = "Scikit"
$ws.Range("B78").Value = "Scikit-learn"
$ws.Range("A79").Value = "sklearn"
$ws.Range("B79").Value = "Scikit-learn"
$ws.Range("A80").Value = "SQL"
$ws.Range("B80").Value = "Structured Query Language (SQL)"
$ws.Range("A81").Value = "Visualizations"
$ws.Range("B81").Value = "Visualization"

# --- 3. Re-apply the handful of cell formats (borders/font) that belong on
#        this refreshed block. Copy *formats only* from cells that already carry
#        the desired style so no new style entries are created in styles.xml. ---
$ws.Range("A48").Copy()
$ws.Range("A34,A41,A61,B62,B66").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A55").Copy()
$ws.Range("A40").PasteSpecial(-4122)

$ws.Range("B56").Copy()
$ws.Range("B42").PasteSpecial(-4122)

$ws.Range("B74").Copy()
$ws.Range("B59").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 4. The donor cells above inherited a style from the pre-edit layout that no
#        longer belongs on this refreshed block; clear it back to the default. ---
$ws.Range("A48").ClearFormats()
$ws.Range("A55").ClearFormats()
$ws.Range("B56").ClearFormats()
$ws.Range("A69").ClearFormats()
$ws.Range("B74").ClearFormats()

# --- 5. Resize Table13 / its AutoFilter to cover the new A1:B81 extent ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B81"))

# --- 6. Match the saved selection / active cell from the edit ---
$ws.Range("A15").Select()

